# Inserts a new weekly price record as row 7 in the daily-logic subset
# sheet, pushing every existing row from 7 downward by one (Excel's
# normal "insert row" shift behavior). The new row carries the same
# market/category/unit metadata as every other row in this sheet and a
# fresh date + price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7..130 down to 8..131, inheriting formatting from the row above
# (matches native Excel "Insert" behavior, including the date style on column D).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new observation.
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44515
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 100112037
$ws.Cells.Item(7, 7).Value = "Cebollín"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 3000
$ws.Cells.Item(7, 11).Value = 900
$ws.Cells.Item(7, 12).Value = 1000
$ws.Cells.Item(7, 13).Value = 950
$ws.Cells.Item(7, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(7, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(7, 16).Value = 158
$ws.Cells.Item(7, 17).Value = 6
$ws.Cells.Item(7, 18).Value = "Hortaliza"
